$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new rows 1021-1098 (UDS Forms B7 + B8 translation-dictionary entries) ---
# row 1021
$ws.Cells.Item(1021, 1).Value = 'b7ptid'
$ws.Cells.Item(1021, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1021, 3).Value = 'b7ptid'
$ws.Cells.Item(1021, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1021, 7).Value = 'Yes'
$ws.Cells.Item(1021, 9).Value = 'LDNM'
# row 1022
$ws.Cells.Item(1022, 1).Value = 'b7pkt_type'
$ws.Cells.Item(1022, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1022, 3).Value = 'b7pkt_type'
$ws.Cells.Item(1022, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1022, 7).Value = 'Yes'
$ws.Cells.Item(1022, 9).Value = 'LDNM'
# row 1023
$ws.Cells.Item(1023, 1).Value = 'b7_formver'
$ws.Cells.Item(1023, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1023, 3).Value = 'b7_formver'
$ws.Cells.Item(1023, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1023, 7).Value = 'Yes'
$ws.Cells.Item(1023, 9).Value = 'LDNM'
# row 1024
$ws.Cells.Item(1024, 1).Value = 'b7form_date'
$ws.Cells.Item(1024, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1024, 3).Value = 'b7form_date'
$ws.Cells.Item(1024, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1024, 7).Value = 'Yes'
$ws.Cells.Item(1024, 9).Value = 'LDNM'
# row 1025
$ws.Cells.Item(1025, 1).Value = 'b7_visit_month'
$ws.Cells.Item(1025, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1025, 3).Value = 'b7_visit_month'
$ws.Cells.Item(1025, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1025, 7).Value = 'Yes'
$ws.Cells.Item(1025, 9).Value = 'LDNM'
# row 1026
$ws.Cells.Item(1026, 1).Value = 'b7_visit_day'
$ws.Cells.Item(1026, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1026, 3).Value = 'b7_visit_day'
$ws.Cells.Item(1026, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1026, 7).Value = 'Yes'
$ws.Cells.Item(1026, 9).Value = 'LDNM'
# row 1027
$ws.Cells.Item(1027, 1).Value = 'b7_visit_yr'
$ws.Cells.Item(1027, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1027, 3).Value = 'b7_visit_yr'
$ws.Cells.Item(1027, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1027, 7).Value = 'Yes'
$ws.Cells.Item(1027, 9).Value = 'LDNM'
# row 1028
$ws.Cells.Item(1028, 1).Value = 'b7visit_num'
$ws.Cells.Item(1028, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1028, 3).Value = 'b7visit_num'
$ws.Cells.Item(1028, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1028, 7).Value = 'Yes'
$ws.Cells.Item(1028, 9).Value = 'LDNM'
# row 1029
$ws.Cells.Item(1029, 1).Value = 'b7_examiner'
$ws.Cells.Item(1029, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1029, 3).Value = 'b7_examiner'
$ws.Cells.Item(1029, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1029, 7).Value = 'Yes'
$ws.Cells.Item(1029, 9).Value = 'LDNM'
# row 1030
$ws.Cells.Item(1030, 1).Value = 'bills'
$ws.Cells.Item(1030, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1030, 3).Value = 'bills'
$ws.Cells.Item(1030, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1030, 5).Value = 'bills'
$ws.Cells.Item(1030, 5).Style = 'Normal'
$ws.Cells.Item(1030, 6).Value = 'ivp_b7'
$ws.Cells.Item(1030, 6).Style = 'Normal'
$ws.Cells.Item(1030, 7).Value = 'Yes'
$ws.Cells.Item(1030, 9).Value = 'LDNM'
# row 1031
$ws.Cells.Item(1031, 1).Value = 'taxes'
$ws.Cells.Item(1031, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1031, 3).Value = 'taxes'
$ws.Cells.Item(1031, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1031, 5).Value = 'taxes'
$ws.Cells.Item(1031, 5).Style = 'Normal'
$ws.Cells.Item(1031, 6).Value = 'ivp_b7'
$ws.Cells.Item(1031, 6).Style = 'Normal'
$ws.Cells.Item(1031, 7).Value = 'Yes'
$ws.Cells.Item(1031, 9).Value = 'LDNM'
# row 1032
$ws.Cells.Item(1032, 1).Value = 'shopping'
$ws.Cells.Item(1032, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1032, 3).Value = 'shopping'
$ws.Cells.Item(1032, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1032, 5).Value = 'shopping'
$ws.Cells.Item(1032, 5).Style = 'Normal'
$ws.Cells.Item(1032, 6).Value = 'ivp_b7'
$ws.Cells.Item(1032, 6).Style = 'Normal'
$ws.Cells.Item(1032, 7).Value = 'Yes'
$ws.Cells.Item(1032, 9).Value = 'LDNM'
# row 1033
$ws.Cells.Item(1033, 1).Value = 'games'
$ws.Cells.Item(1033, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1033, 3).Value = 'games'
$ws.Cells.Item(1033, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1033, 5).Value = 'games'
$ws.Cells.Item(1033, 5).Style = 'Normal'
$ws.Cells.Item(1033, 6).Value = 'ivp_b7'
$ws.Cells.Item(1033, 6).Style = 'Normal'
$ws.Cells.Item(1033, 7).Value = 'Yes'
$ws.Cells.Item(1033, 9).Value = 'LDNM'
# row 1034
$ws.Cells.Item(1034, 1).Value = 'stove'
$ws.Cells.Item(1034, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1034, 3).Value = 'stove'
$ws.Cells.Item(1034, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1034, 5).Value = 'stove'
$ws.Cells.Item(1034, 5).Style = 'Normal'
$ws.Cells.Item(1034, 6).Value = 'ivp_b7'
$ws.Cells.Item(1034, 6).Style = 'Normal'
$ws.Cells.Item(1034, 7).Value = 'Yes'
$ws.Cells.Item(1034, 9).Value = 'LDNM'
# row 1035
$ws.Cells.Item(1035, 1).Value = 'mealprep'
$ws.Cells.Item(1035, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1035, 3).Value = 'mealprep'
$ws.Cells.Item(1035, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1035, 5).Value = 'mealprep'
$ws.Cells.Item(1035, 5).Style = 'Normal'
$ws.Cells.Item(1035, 6).Value = 'ivp_b7'
$ws.Cells.Item(1035, 6).Style = 'Normal'
$ws.Cells.Item(1035, 7).Value = 'Yes'
$ws.Cells.Item(1035, 9).Value = 'LDNM'
# row 1036
$ws.Cells.Item(1036, 1).Value = 'events'
$ws.Cells.Item(1036, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1036, 3).Value = 'events'
$ws.Cells.Item(1036, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1036, 5).Value = 'events'
$ws.Cells.Item(1036, 5).Style = 'Normal'
$ws.Cells.Item(1036, 6).Value = 'ivp_b7'
$ws.Cells.Item(1036, 6).Style = 'Normal'
$ws.Cells.Item(1036, 7).Value = 'Yes'
$ws.Cells.Item(1036, 9).Value = 'LDNM'
# row 1037
$ws.Cells.Item(1037, 1).Value = 'payattn'
$ws.Cells.Item(1037, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1037, 3).Value = 'payattn'
$ws.Cells.Item(1037, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1037, 5).Value = 'payattn'
$ws.Cells.Item(1037, 5).Style = 'Normal'
$ws.Cells.Item(1037, 6).Value = 'ivp_b7'
$ws.Cells.Item(1037, 6).Style = 'Normal'
$ws.Cells.Item(1037, 7).Value = 'Yes'
$ws.Cells.Item(1037, 9).Value = 'LDNM'
# row 1038
$ws.Cells.Item(1038, 1).Value = 'remdates'
$ws.Cells.Item(1038, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1038, 3).Value = 'remdates'
$ws.Cells.Item(1038, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1038, 5).Value = 'remdates'
$ws.Cells.Item(1038, 5).Style = 'Normal'
$ws.Cells.Item(1038, 6).Value = 'ivp_b7'
$ws.Cells.Item(1038, 6).Style = 'Normal'
$ws.Cells.Item(1038, 7).Value = 'Yes'
$ws.Cells.Item(1038, 9).Value = 'LDNM'
# row 1039
$ws.Cells.Item(1039, 1).Value = 'travel'
$ws.Cells.Item(1039, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1039, 3).Value = 'travel'
$ws.Cells.Item(1039, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1039, 5).Value = 'travel'
$ws.Cells.Item(1039, 5).Style = 'Normal'
$ws.Cells.Item(1039, 6).Value = 'ivp_b7'
$ws.Cells.Item(1039, 6).Style = 'Normal'
$ws.Cells.Item(1039, 7).Value = 'Yes'
$ws.Cells.Item(1039, 9).Value = 'LDNM'
# row 1040
$ws.Cells.Item(1040, 1).Value = 'form_b7_functional_assessment_faq_complete'
$ws.Cells.Item(1040, 2).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1040, 3).Value = 'form_b7_functional_assessment_faq_complete'
$ws.Cells.Item(1040, 4).Value = 'form_b7_functional_assessment_faq'
$ws.Cells.Item(1040, 5).Value = 'ivp_b7_complete'
$ws.Cells.Item(1040, 5).Style = 'Normal'
$ws.Cells.Item(1040, 6).Value = 'ivp_b7'
$ws.Cells.Item(1040, 6).Style = 'Normal'
$ws.Cells.Item(1040, 7).Value = 'Yes'
$ws.Cells.Item(1040, 9).Value = 'LDNM'
# row 1041
$ws.Cells.Item(1041, 5).Value = 'fas_score'
$ws.Cells.Item(1041, 5).Style = 'Normal'
$ws.Cells.Item(1041, 6).Value = 'ivp_b7'
$ws.Cells.Item(1041, 6).Style = 'Normal'
$ws.Cells.Item(1041, 7).Value = 'Yes'
$ws.Cells.Item(1041, 9).Value = 'LDNM'
# row 1042
$ws.Cells.Item(1042, 1).Value = 'b8ptid'
$ws.Cells.Item(1042, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1042, 3).Value = 'b8ptid'
$ws.Cells.Item(1042, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1042, 7).Value = 'Yes'
$ws.Cells.Item(1042, 9).Value = 'LDNM'
# row 1043
$ws.Cells.Item(1043, 1).Value = 'b8pkt_type'
$ws.Cells.Item(1043, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1043, 3).Value = 'b8pkt_type'
$ws.Cells.Item(1043, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1043, 7).Value = 'Yes'
$ws.Cells.Item(1043, 9).Value = 'LDNM'
# row 1044
$ws.Cells.Item(1044, 1).Value = 'b8_formver'
$ws.Cells.Item(1044, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1044, 3).Value = 'b8_formver'
$ws.Cells.Item(1044, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1044, 7).Value = 'Yes'
$ws.Cells.Item(1044, 9).Value = 'LDNM'
# row 1045
$ws.Cells.Item(1045, 1).Value = 'b8form_date'
$ws.Cells.Item(1045, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1045, 3).Value = 'b8form_date'
$ws.Cells.Item(1045, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1045, 7).Value = 'Yes'
$ws.Cells.Item(1045, 9).Value = 'LDNM'
# row 1046
$ws.Cells.Item(1046, 1).Value = 'b8_visit_month'
$ws.Cells.Item(1046, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1046, 3).Value = 'b8_visit_month'
$ws.Cells.Item(1046, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1046, 7).Value = 'Yes'
$ws.Cells.Item(1046, 9).Value = 'LDNM'
# row 1047
$ws.Cells.Item(1047, 1).Value = 'b8_visit_day'
$ws.Cells.Item(1047, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1047, 3).Value = 'b8_visit_day'
$ws.Cells.Item(1047, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1047, 7).Value = 'Yes'
$ws.Cells.Item(1047, 9).Value = 'LDNM'
# row 1048
$ws.Cells.Item(1048, 1).Value = 'b8_visit_yr'
$ws.Cells.Item(1048, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1048, 3).Value = 'b8_visit_yr'
$ws.Cells.Item(1048, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1048, 7).Value = 'Yes'
$ws.Cells.Item(1048, 9).Value = 'LDNM'
# row 1049
$ws.Cells.Item(1049, 1).Value = 'b8visit_num'
$ws.Cells.Item(1049, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1049, 3).Value = 'b8visit_num'
$ws.Cells.Item(1049, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1049, 7).Value = 'Yes'
$ws.Cells.Item(1049, 9).Value = 'LDNM'
# row 1050
$ws.Cells.Item(1050, 1).Value = 'b8_examiner'
$ws.Cells.Item(1050, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1050, 3).Value = 'b8_examiner'
$ws.Cells.Item(1050, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1050, 7).Value = 'Yes'
$ws.Cells.Item(1050, 9).Value = 'LDNM'
# row 1051
$ws.Cells.Item(1051, 1).Value = 'normal'
$ws.Cells.Item(1051, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1051, 7).Value = 'Yes'
$ws.Cells.Item(1051, 9).Value = 'LDNM'
# row 1052
$ws.Cells.Item(1052, 1).Value = 'focldef'
$ws.Cells.Item(1052, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1052, 7).Value = 'Yes'
$ws.Cells.Item(1052, 9).Value = 'LDNM'
# row 1053
$ws.Cells.Item(1053, 1).Value = 'gaitdis'
$ws.Cells.Item(1053, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1053, 7).Value = 'Yes'
$ws.Cells.Item(1053, 9).Value = 'LDNM'
# row 1054
$ws.Cells.Item(1054, 1).Value = 'eyemove'
$ws.Cells.Item(1054, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1054, 7).Value = 'Yes'
$ws.Cells.Item(1054, 9).Value = 'LDNM'
# row 1055
$ws.Cells.Item(1055, 1).Value = 'form_b8_physical_neurological_exam_findings_complete'
$ws.Cells.Item(1055, 2).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1055, 3).Value = 'form_b8_physical_neurological_exam_findings_complete'
$ws.Cells.Item(1055, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1055, 5).Value = 'ivp_b8_complete'
$ws.Cells.Item(1055, 5).Style = 'Normal'
$ws.Cells.Item(1055, 6).Value = 'ivp_b8'
$ws.Cells.Item(1055, 6).Style = 'Normal'
$ws.Cells.Item(1055, 7).Value = 'Yes'
$ws.Cells.Item(1055, 9).Value = 'LDNM'
# row 1056
$ws.Cells.Item(1056, 3).Value = 'normexam'
$ws.Cells.Item(1056, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1056, 5).Value = 'normexam'
$ws.Cells.Item(1056, 5).Style = 'Normal'
$ws.Cells.Item(1056, 6).Value = 'ivp_b8'
$ws.Cells.Item(1056, 6).Style = 'Normal'
$ws.Cells.Item(1056, 7).Value = 'Yes'
$ws.Cells.Item(1056, 9).Value = 'LDNM'
# row 1057
$ws.Cells.Item(1057, 3).Value = 'parksign'
$ws.Cells.Item(1057, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1057, 5).Value = 'parksign'
$ws.Cells.Item(1057, 5).Style = 'Normal'
$ws.Cells.Item(1057, 6).Value = 'ivp_b8'
$ws.Cells.Item(1057, 6).Style = 'Normal'
$ws.Cells.Item(1057, 7).Value = 'Yes'
$ws.Cells.Item(1057, 9).Value = 'LDNM'
# row 1058
$ws.Cells.Item(1058, 3).Value = 'resttrl'
$ws.Cells.Item(1058, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1058, 5).Value = 'resttrl'
$ws.Cells.Item(1058, 5).Style = 'Normal'
$ws.Cells.Item(1058, 6).Value = 'ivp_b8'
$ws.Cells.Item(1058, 6).Style = 'Normal'
$ws.Cells.Item(1058, 7).Value = 'Yes'
$ws.Cells.Item(1058, 9).Value = 'LDNM'
# row 1059
$ws.Cells.Item(1059, 3).Value = 'resttrr'
$ws.Cells.Item(1059, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1059, 5).Value = 'resttrr'
$ws.Cells.Item(1059, 5).Style = 'Normal'
$ws.Cells.Item(1059, 6).Value = 'ivp_b8'
$ws.Cells.Item(1059, 6).Style = 'Normal'
$ws.Cells.Item(1059, 7).Value = 'Yes'
$ws.Cells.Item(1059, 9).Value = 'LDNM'
# row 1060
$ws.Cells.Item(1060, 3).Value = 'slowingl'
$ws.Cells.Item(1060, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1060, 5).Value = 'slowingl'
$ws.Cells.Item(1060, 5).Style = 'Normal'
$ws.Cells.Item(1060, 6).Value = 'ivp_b8'
$ws.Cells.Item(1060, 6).Style = 'Normal'
$ws.Cells.Item(1060, 7).Value = 'Yes'
$ws.Cells.Item(1060, 9).Value = 'LDNM'
# row 1061
$ws.Cells.Item(1061, 3).Value = 'slowingr'
$ws.Cells.Item(1061, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1061, 5).Value = 'slowingr'
$ws.Cells.Item(1061, 5).Style = 'Normal'
$ws.Cells.Item(1061, 6).Value = 'ivp_b8'
$ws.Cells.Item(1061, 6).Style = 'Normal'
$ws.Cells.Item(1061, 7).Value = 'Yes'
$ws.Cells.Item(1061, 9).Value = 'LDNM'
# row 1062
$ws.Cells.Item(1062, 3).Value = 'rigidl'
$ws.Cells.Item(1062, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1062, 5).Value = 'rigidl'
$ws.Cells.Item(1062, 5).Style = 'Normal'
$ws.Cells.Item(1062, 6).Value = 'ivp_b8'
$ws.Cells.Item(1062, 6).Style = 'Normal'
$ws.Cells.Item(1062, 7).Value = 'Yes'
$ws.Cells.Item(1062, 9).Value = 'LDNM'
# row 1063
$ws.Cells.Item(1063, 3).Value = 'rigidr'
$ws.Cells.Item(1063, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1063, 5).Value = 'rigidr'
$ws.Cells.Item(1063, 5).Style = 'Normal'
$ws.Cells.Item(1063, 6).Value = 'ivp_b8'
$ws.Cells.Item(1063, 6).Style = 'Normal'
$ws.Cells.Item(1063, 7).Value = 'Yes'
$ws.Cells.Item(1063, 9).Value = 'LDNM'
# row 1064
$ws.Cells.Item(1064, 3).Value = 'brady'
$ws.Cells.Item(1064, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1064, 5).Value = 'brady'
$ws.Cells.Item(1064, 5).Style = 'Normal'
$ws.Cells.Item(1064, 6).Value = 'ivp_b8'
$ws.Cells.Item(1064, 6).Style = 'Normal'
$ws.Cells.Item(1064, 7).Value = 'Yes'
$ws.Cells.Item(1064, 9).Value = 'LDNM'
# row 1065
$ws.Cells.Item(1065, 3).Value = 'parkgait'
$ws.Cells.Item(1065, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1065, 5).Value = 'parkgait'
$ws.Cells.Item(1065, 5).Style = 'Normal'
$ws.Cells.Item(1065, 6).Value = 'ivp_b8'
$ws.Cells.Item(1065, 6).Style = 'Normal'
$ws.Cells.Item(1065, 7).Value = 'Yes'
$ws.Cells.Item(1065, 9).Value = 'LDNM'
# row 1066
$ws.Cells.Item(1066, 3).Value = 'postinst'
$ws.Cells.Item(1066, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1066, 5).Value = 'postinst'
$ws.Cells.Item(1066, 5).Style = 'Normal'
$ws.Cells.Item(1066, 6).Value = 'ivp_b8'
$ws.Cells.Item(1066, 6).Style = 'Normal'
$ws.Cells.Item(1066, 7).Value = 'Yes'
$ws.Cells.Item(1066, 9).Value = 'LDNM'
# row 1067
$ws.Cells.Item(1067, 3).Value = 'cvdsigns'
$ws.Cells.Item(1067, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1067, 5).Value = 'cvdsigns'
$ws.Cells.Item(1067, 5).Style = 'Normal'
$ws.Cells.Item(1067, 6).Value = 'ivp_b8'
$ws.Cells.Item(1067, 6).Style = 'Normal'
$ws.Cells.Item(1067, 7).Value = 'Yes'
$ws.Cells.Item(1067, 9).Value = 'LDNM'
# row 1068
$ws.Cells.Item(1068, 3).Value = 'cortdef'
$ws.Cells.Item(1068, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1068, 5).Value = 'cortdef'
$ws.Cells.Item(1068, 5).Style = 'Normal'
$ws.Cells.Item(1068, 6).Value = 'ivp_b8'
$ws.Cells.Item(1068, 6).Style = 'Normal'
$ws.Cells.Item(1068, 7).Value = 'Yes'
$ws.Cells.Item(1068, 9).Value = 'LDNM'
# row 1069
$ws.Cells.Item(1069, 3).Value = 'sivdfind'
$ws.Cells.Item(1069, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1069, 5).Value = 'sivdfind'
$ws.Cells.Item(1069, 5).Style = 'Normal'
$ws.Cells.Item(1069, 6).Value = 'ivp_b8'
$ws.Cells.Item(1069, 6).Style = 'Normal'
$ws.Cells.Item(1069, 7).Value = 'Yes'
$ws.Cells.Item(1069, 9).Value = 'LDNM'
# row 1070
$ws.Cells.Item(1070, 3).Value = 'cvdmotl'
$ws.Cells.Item(1070, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1070, 5).Value = 'cvdmotl'
$ws.Cells.Item(1070, 5).Style = 'Normal'
$ws.Cells.Item(1070, 6).Value = 'ivp_b8'
$ws.Cells.Item(1070, 6).Style = 'Normal'
$ws.Cells.Item(1070, 7).Value = 'Yes'
$ws.Cells.Item(1070, 9).Value = 'LDNM'
# row 1071
$ws.Cells.Item(1071, 3).Value = 'cvdmotr'
$ws.Cells.Item(1071, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1071, 5).Value = 'cvdmotr'
$ws.Cells.Item(1071, 5).Style = 'Normal'
$ws.Cells.Item(1071, 6).Value = 'ivp_b8'
$ws.Cells.Item(1071, 6).Style = 'Normal'
$ws.Cells.Item(1071, 7).Value = 'Yes'
$ws.Cells.Item(1071, 9).Value = 'LDNM'
# row 1072
$ws.Cells.Item(1072, 3).Value = 'cortvisl'
$ws.Cells.Item(1072, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1072, 5).Value = 'cortvisl'
$ws.Cells.Item(1072, 5).Style = 'Normal'
$ws.Cells.Item(1072, 6).Value = 'ivp_b8'
$ws.Cells.Item(1072, 6).Style = 'Normal'
$ws.Cells.Item(1072, 7).Value = 'Yes'
$ws.Cells.Item(1072, 9).Value = 'LDNM'
# row 1073
$ws.Cells.Item(1073, 3).Value = 'cortvisr'
$ws.Cells.Item(1073, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1073, 5).Value = 'cortvisr'
$ws.Cells.Item(1073, 5).Style = 'Normal'
$ws.Cells.Item(1073, 6).Value = 'ivp_b8'
$ws.Cells.Item(1073, 6).Style = 'Normal'
$ws.Cells.Item(1073, 7).Value = 'Yes'
$ws.Cells.Item(1073, 9).Value = 'LDNM'
# row 1074
$ws.Cells.Item(1074, 3).Value = 'somatl'
$ws.Cells.Item(1074, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1074, 5).Value = 'somatl'
$ws.Cells.Item(1074, 5).Style = 'Normal'
$ws.Cells.Item(1074, 6).Value = 'ivp_b8'
$ws.Cells.Item(1074, 6).Style = 'Normal'
$ws.Cells.Item(1074, 7).Value = 'Yes'
$ws.Cells.Item(1074, 9).Value = 'LDNM'
# row 1075
$ws.Cells.Item(1075, 3).Value = 'somatr'
$ws.Cells.Item(1075, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1075, 5).Value = 'somatr'
$ws.Cells.Item(1075, 5).Style = 'Normal'
$ws.Cells.Item(1075, 6).Value = 'ivp_b8'
$ws.Cells.Item(1075, 6).Style = 'Normal'
$ws.Cells.Item(1075, 7).Value = 'Yes'
$ws.Cells.Item(1075, 9).Value = 'LDNM'
# row 1076
$ws.Cells.Item(1076, 3).Value = 'postcort'
$ws.Cells.Item(1076, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1076, 5).Value = 'postcort'
$ws.Cells.Item(1076, 5).Style = 'Normal'
$ws.Cells.Item(1076, 6).Value = 'ivp_b8'
$ws.Cells.Item(1076, 6).Style = 'Normal'
$ws.Cells.Item(1076, 7).Value = 'Yes'
$ws.Cells.Item(1076, 9).Value = 'LDNM'
# row 1077
$ws.Cells.Item(1077, 3).Value = 'pspcbs'
$ws.Cells.Item(1077, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1077, 5).Value = 'pspcbs'
$ws.Cells.Item(1077, 5).Style = 'Normal'
$ws.Cells.Item(1077, 6).Value = 'ivp_b8'
$ws.Cells.Item(1077, 6).Style = 'Normal'
$ws.Cells.Item(1077, 7).Value = 'Yes'
$ws.Cells.Item(1077, 9).Value = 'LDNM'
# row 1078
$ws.Cells.Item(1078, 3).Value = 'eyepsp'
$ws.Cells.Item(1078, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1078, 5).Value = 'eyepsp'
$ws.Cells.Item(1078, 5).Style = 'Normal'
$ws.Cells.Item(1078, 6).Value = 'ivp_b8'
$ws.Cells.Item(1078, 6).Style = 'Normal'
$ws.Cells.Item(1078, 7).Value = 'Yes'
$ws.Cells.Item(1078, 9).Value = 'LDNM'
# row 1079
$ws.Cells.Item(1079, 3).Value = 'dyspsp'
$ws.Cells.Item(1079, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1079, 5).Value = 'dyspsp'
$ws.Cells.Item(1079, 5).Style = 'Normal'
$ws.Cells.Item(1079, 6).Value = 'ivp_b8'
$ws.Cells.Item(1079, 6).Style = 'Normal'
$ws.Cells.Item(1079, 7).Value = 'Yes'
$ws.Cells.Item(1079, 9).Value = 'LDNM'
# row 1080
$ws.Cells.Item(1080, 3).Value = 'axialpsp'
$ws.Cells.Item(1080, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1080, 5).Value = 'axialpsp'
$ws.Cells.Item(1080, 5).Style = 'Normal'
$ws.Cells.Item(1080, 6).Value = 'ivp_b8'
$ws.Cells.Item(1080, 6).Style = 'Normal'
$ws.Cells.Item(1080, 7).Value = 'Yes'
$ws.Cells.Item(1080, 9).Value = 'LDNM'
# row 1081
$ws.Cells.Item(1081, 3).Value = 'gaitpsp'
$ws.Cells.Item(1081, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1081, 5).Value = 'gaitpsp'
$ws.Cells.Item(1081, 5).Style = 'Normal'
$ws.Cells.Item(1081, 6).Value = 'ivp_b8'
$ws.Cells.Item(1081, 6).Style = 'Normal'
$ws.Cells.Item(1081, 7).Value = 'Yes'
$ws.Cells.Item(1081, 9).Value = 'LDNM'
# row 1082
$ws.Cells.Item(1082, 3).Value = 'apraxsp'
$ws.Cells.Item(1082, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1082, 5).Value = 'apraxsp'
$ws.Cells.Item(1082, 5).Style = 'Normal'
$ws.Cells.Item(1082, 6).Value = 'ivp_b8'
$ws.Cells.Item(1082, 6).Style = 'Normal'
$ws.Cells.Item(1082, 7).Value = 'Yes'
$ws.Cells.Item(1082, 9).Value = 'LDNM'
# row 1083
$ws.Cells.Item(1083, 3).Value = 'apraxl'
$ws.Cells.Item(1083, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1083, 5).Value = 'apraxl'
$ws.Cells.Item(1083, 5).Style = 'Normal'
$ws.Cells.Item(1083, 6).Value = 'ivp_b8'
$ws.Cells.Item(1083, 6).Style = 'Normal'
$ws.Cells.Item(1083, 7).Value = 'Yes'
$ws.Cells.Item(1083, 9).Value = 'LDNM'
# row 1084
$ws.Cells.Item(1084, 3).Value = 'apraxr'
$ws.Cells.Item(1084, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1084, 5).Value = 'apraxr'
$ws.Cells.Item(1084, 5).Style = 'Normal'
$ws.Cells.Item(1084, 6).Value = 'ivp_b8'
$ws.Cells.Item(1084, 6).Style = 'Normal'
$ws.Cells.Item(1084, 7).Value = 'Yes'
$ws.Cells.Item(1084, 9).Value = 'LDNM'
# row 1085
$ws.Cells.Item(1085, 3).Value = 'cortsenl'
$ws.Cells.Item(1085, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1085, 5).Value = 'cortsenl'
$ws.Cells.Item(1085, 5).Style = 'Normal'
$ws.Cells.Item(1085, 6).Value = 'ivp_b8'
$ws.Cells.Item(1085, 6).Style = 'Normal'
$ws.Cells.Item(1085, 7).Value = 'Yes'
$ws.Cells.Item(1085, 9).Value = 'LDNM'
# row 1086
$ws.Cells.Item(1086, 3).Value = 'cortsenr'
$ws.Cells.Item(1086, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1086, 5).Value = 'cortsenr'
$ws.Cells.Item(1086, 5).Style = 'Normal'
$ws.Cells.Item(1086, 6).Value = 'ivp_b8'
$ws.Cells.Item(1086, 6).Style = 'Normal'
$ws.Cells.Item(1086, 7).Value = 'Yes'
$ws.Cells.Item(1086, 9).Value = 'LDNM'
# row 1087
$ws.Cells.Item(1087, 3).Value = 'ataxl'
$ws.Cells.Item(1087, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1087, 5).Value = 'ataxl'
$ws.Cells.Item(1087, 5).Style = 'Normal'
$ws.Cells.Item(1087, 6).Value = 'ivp_b8'
$ws.Cells.Item(1087, 6).Style = 'Normal'
$ws.Cells.Item(1087, 7).Value = 'Yes'
$ws.Cells.Item(1087, 9).Value = 'LDNM'
# row 1088
$ws.Cells.Item(1088, 3).Value = 'ataxr'
$ws.Cells.Item(1088, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1088, 5).Value = 'ataxr'
$ws.Cells.Item(1088, 5).Style = 'Normal'
$ws.Cells.Item(1088, 6).Value = 'ivp_b8'
$ws.Cells.Item(1088, 6).Style = 'Normal'
$ws.Cells.Item(1088, 7).Value = 'Yes'
$ws.Cells.Item(1088, 9).Value = 'LDNM'
# row 1089
$ws.Cells.Item(1089, 3).Value = 'akuebknk'
$ws.Cells.Item(1089, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1089, 5).Value = 'alienlml'
$ws.Cells.Item(1089, 5).Style = 'Normal'
$ws.Cells.Item(1089, 6).Value = 'ivp_b8'
$ws.Cells.Item(1089, 6).Style = 'Normal'
$ws.Cells.Item(1089, 7).Value = 'Yes'
$ws.Cells.Item(1089, 9).Value = 'LDNM'
# row 1090
$ws.Cells.Item(1090, 3).Value = 'alienlmr'
$ws.Cells.Item(1090, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1090, 5).Value = 'alienlmr'
$ws.Cells.Item(1090, 5).Style = 'Normal'
$ws.Cells.Item(1090, 6).Value = 'ivp_b8'
$ws.Cells.Item(1090, 6).Style = 'Normal'
$ws.Cells.Item(1090, 7).Value = 'Yes'
$ws.Cells.Item(1090, 9).Value = 'LDNM'
# row 1091
$ws.Cells.Item(1091, 3).Value = 'dystonl'
$ws.Cells.Item(1091, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1091, 5).Value = 'dystonl'
$ws.Cells.Item(1091, 5).Style = 'Normal'
$ws.Cells.Item(1091, 6).Value = 'ivp_b8'
$ws.Cells.Item(1091, 6).Style = 'Normal'
$ws.Cells.Item(1091, 7).Value = 'Yes'
$ws.Cells.Item(1091, 9).Value = 'LDNM'
# row 1092
$ws.Cells.Item(1092, 3).Value = 'dystonr'
$ws.Cells.Item(1092, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1092, 5).Value = 'dystonr'
$ws.Cells.Item(1092, 5).Style = 'Normal'
$ws.Cells.Item(1092, 6).Value = 'ivp_b8'
$ws.Cells.Item(1092, 6).Style = 'Normal'
$ws.Cells.Item(1092, 7).Value = 'Yes'
$ws.Cells.Item(1092, 9).Value = 'LDNM'
# row 1093
$ws.Cells.Item(1093, 3).Value = 'myocllt'
$ws.Cells.Item(1093, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1093, 5).Value = 'myocllt'
$ws.Cells.Item(1093, 5).Style = 'Normal'
$ws.Cells.Item(1093, 6).Value = 'ivp_b8'
$ws.Cells.Item(1093, 6).Style = 'Normal'
$ws.Cells.Item(1093, 7).Value = 'Yes'
$ws.Cells.Item(1093, 9).Value = 'LDNM'
# row 1094
$ws.Cells.Item(1094, 3).Value = 'myoclrt'
$ws.Cells.Item(1094, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1094, 5).Value = 'myoclrt'
$ws.Cells.Item(1094, 5).Style = 'Normal'
$ws.Cells.Item(1094, 6).Value = 'ivp_b8'
$ws.Cells.Item(1094, 6).Style = 'Normal'
$ws.Cells.Item(1094, 7).Value = 'Yes'
$ws.Cells.Item(1094, 9).Value = 'LDNM'
# row 1095
$ws.Cells.Item(1095, 3).Value = 'alsfind'
$ws.Cells.Item(1095, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1095, 5).Value = 'alsfind'
$ws.Cells.Item(1095, 5).Style = 'Normal'
$ws.Cells.Item(1095, 6).Value = 'ivp_b8'
$ws.Cells.Item(1095, 6).Style = 'Normal'
$ws.Cells.Item(1095, 7).Value = 'Yes'
$ws.Cells.Item(1095, 9).Value = 'LDNM'
# row 1096
$ws.Cells.Item(1096, 3).Value = 'gaitnph'
$ws.Cells.Item(1096, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1096, 5).Value = 'gaitnph'
$ws.Cells.Item(1096, 5).Style = 'Normal'
$ws.Cells.Item(1096, 6).Value = 'ivp_b8'
$ws.Cells.Item(1096, 6).Style = 'Normal'
$ws.Cells.Item(1096, 7).Value = 'Yes'
$ws.Cells.Item(1096, 9).Value = 'LDNM'
# row 1097
$ws.Cells.Item(1097, 3).Value = 'othneur'
$ws.Cells.Item(1097, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1097, 5).Value = 'othneur'
$ws.Cells.Item(1097, 5).Style = 'Normal'
$ws.Cells.Item(1097, 6).Value = 'ivp_b8'
$ws.Cells.Item(1097, 6).Style = 'Normal'
$ws.Cells.Item(1097, 7).Value = 'Yes'
$ws.Cells.Item(1097, 9).Value = 'LDNM'
# row 1098
$ws.Cells.Item(1098, 3).Value = 'othneurx'
$ws.Cells.Item(1098, 4).Value = 'form_b8_physical_neurological_exam_findings'
$ws.Cells.Item(1098, 5).Value = 'othneurx'
$ws.Cells.Item(1098, 5).Style = 'Normal'
$ws.Cells.Item(1098, 6).Value = 'ivp_b8'
$ws.Cells.Item(1098, 6).Style = 'Normal'
$ws.Cells.Item(1098, 7).Value = 'Yes'
$ws.Cells.Item(1098, 9).Value = 'LDNM'

# --- Column width adjustments: columns C and D mirror A and B's bestfit widths ---
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# --- Update selection to reflect the new bottom of the data range ---
$ws.Range("A1099").Select()
